# Undo Jason's overwrite of class materials starter code/slides.
#
# The only author-intended content change that is reachable through the
# PowerPoint object model is on slide 1: the placeholder password text
# "errorsasobjects" (which PowerPoint had flagged as a possible
# misspelling) is replaced with a blank-style run of underscores
# ("__________") inside the "Today's Attendance password" textbox,
# while keeping the existing run formatting (size, highlight, dirty flag).

$p = $ppt.ActivePresentation

$targetOld = "errorsasobjects"
$targetNew = "__________"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            $full = $tr.Text
            if ($full.IndexOf($targetOld) -ge 0) {
                [void]$tr.Replace($targetOld, $targetNew)
            }
        }
    }
}
